$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing work-day entry for row 18 (date 2020-09-10),
# matching the style of the surrounding rows (date + time formats).
$ws.Range("A18").Value = 44084
$ws.Range("A18").NumberFormat = "d-mmm"

$ws.Range("B18").Value = 0.70833333333333337
$ws.Range("B18").NumberFormat = "h:mm"

$ws.Range("C18").Value = 0.91666666666666663
$ws.Range("C18").NumberFormat = "h:mm"

# Move the active selection to C9, matching the saved cursor position.
$ws.Range("C9").Select()
